# Apply updated cryptocurrency data (prices, 1h volume %, and
# coin/link reshuffle for rows 15-23) per the "Updated symbol list"
# commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '306.79'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.64%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-7.92%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.091'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.53%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07689'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-6.61%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.246'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.63%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.625'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-8.85%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9154'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.41%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1016'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-9.61%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1749'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-7.82%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09372'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.64%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04437'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.09%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1057'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.11%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001258'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-3.47%'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04143'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.56%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005826'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.25%'
$ws.Range("B17").Value = 'UpBots'
$ws.Range("C17").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.007491'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2,414.23%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.362'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.02%'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.433'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-4.28%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3304'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.84%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.030'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-5.55%'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1348'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.99%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.2811'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '10.28%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001199'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.97%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004090'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.44%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '6.70%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02427'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-9.60%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05178'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-7.52%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.13%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-6.24%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007171'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '9.50%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001950'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-7.84%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008369'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '9.17%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3055'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-12.29%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006415'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-5.32%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-26.81%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004372'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '42.45%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
